# "added 4wk low sales check"
# Update the per-week forecast figures (MyForecast, Inventory Coverage,
# Stockout Risk, Seasonality Index) on the "Forecast Comparison" sheet,
# and the derived totals on the "Summary" sheet.

$wb = $excel.ActiveWorkbook
$wsForecast = $wb.Worksheets.Item("Forecast Comparison")
$wsSummary  = $wb.Worksheets.Item("Summary")

# Row => MyForecast (D), Inventory Coverage (H), Stockout Risk (I), Seasonality Index (L)
$rows = @(
    @{ Row = 2;  D = 80;  H = 7.62; I = $null;  L = 1.13 },
    @{ Row = 3;  D = 84;  H = 6.31; I = $null;  L = 0.86 },
    @{ Row = 4;  D = 88;  H = 5.07; I = $null;  L = 1 },
    @{ Row = 5;  D = 88;  H = 4.07; I = $null;  L = 1.05 },
    @{ Row = 6;  D = 86;  H = 3.14; I = $null;  L = 1.14 },
    @{ Row = 7;  D = 87;  H = 2.11; I = $null;  L = 1.13 },
    @{ Row = 8;  D = 93;  H = 1.04; I = $null;  L = 1.01 },
    @{ Row = 9;  D = 100; H = 0.04; I = "High"; L = 1.07 },
    @{ Row = 10; D = 103; H = $null; I = $null; L = 0.89 },
    @{ Row = 11; D = 101; H = $null; I = $null; L = 0.91 },
    @{ Row = 12; D = 100; H = $null; I = $null; L = 1.03 },
    @{ Row = 13; D = 105; H = $null; I = $null; L = 1.05 },
    @{ Row = 14; D = 112; H = $null; I = $null; L = 0.87 },
    @{ Row = 15; D = 117; H = $null; I = $null; L = 1.01 },
    @{ Row = 16; D = 116; H = $null; I = $null; L = 1.11 },
    @{ Row = 17; D = 114; H = $null; I = $null; L = 0.9399999999999999 }
)

foreach ($r in $rows) {
    $wsForecast.Range("D$($r.Row)").Value = $r.D
    if ($null -ne $r.H) {
        $wsForecast.Range("H$($r.Row)").Value = $r.H
    }
    if ($null -ne $r.I) {
        $wsForecast.Range("I$($r.Row)").Value = $r.I
    }
    $wsForecast.Range("L$($r.Row)").Value = $r.L
}

# Update the derived summary totals. These cells hold their numbers as
# text (matching the sheet's existing inline-string formatting), so a
# leading apostrophe is used to force text entry instead of letting the
# numeric-looking value get auto-converted to a number. ClearFormats()
# then drops the incidental "quote prefix" formatting flag that the
# apostrophe entry leaves behind, restoring the cell's original (default)
# styling.
$wsSummary.Range("B9").Value2  = "'1574"
$wsSummary.Range("B9").ClearFormats()
$wsSummary.Range("B10").Value2 = "'706"
$wsSummary.Range("B10").ClearFormats()
$wsSummary.Range("B11").Value2 = "'340"
$wsSummary.Range("B11").ClearFormats()
$wsSummary.Range("B12").Value2 = "'117"
$wsSummary.Range("B12").ClearFormats()
$wsSummary.Range("B14").Value2 = "'80"
$wsSummary.Range("B14").ClearFormats()
